$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '69.616.65'
$ws.Cells.Item(2, 5).Value = '  -1.74%  '
$ws.Cells.Item(3, 4).Value = '3.511.26'
$ws.Cells.Item(3, 5).Value = '  -1.32%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).Value = '''616.47'
$ws.Cells.Item(5, 5).Value = '  +5.76%  '
$ws.Cells.Item(6, 4).Value = '''191.87'
$ws.Cells.Item(6, 5).Value = '  +1.64%  '
$ws.Cells.Item(7, 5).Value = '  +0.26%  '
$ws.Cells.Item(8, 5).Value = '  -0.04%  '
$ws.Cells.Item(9, 4).Value = '''0.211'
$ws.Cells.Item(9, 5).Value = '  -3.16%  '
$ws.Cells.Item(10, 4).Value = '''0.654'
$ws.Cells.Item(10, 5).Value = '  +0.61%  '
$ws.Cells.Item(11, 4).Value = '''53.46'
$ws.Cells.Item(11, 5).Value = '  -1.92%  '
$ws.Cells.Item(12, 5).Value = '  -3.27%  '
$ws.Cells.Item(13, 4).Value = '''9.55'
$ws.Cells.Item(13, 5).Value = '  +0.38%  '
$ws.Cells.Item(14, 4).Value = '4.075.90'
$ws.Cells.Item(14, 5).Value = '  -1.08%  '
$ws.Cells.Item(15, 4).Value = '''610.79'
$ws.Cells.Item(15, 5).Value = '  +6.18%  '
$ws.Cells.Item(16, 4).Value = '69.742.49'
$ws.Cells.Item(16, 5).Value = '  -1.55%  '
$ws.Cells.Item(17, 4).Value = '''18.99'
$ws.Cells.Item(17, 5).Value = '  -1.01%  '
$ws.Cells.Item(18, 4).Value = '''12.57'
$ws.Cells.Item(18, 5).Value = '  -1.49%  '
$ws.Cells.Item(19, 4).Value = '3.516.54'
$ws.Cells.Item(19, 5).Value = '  -1.18%  '
$ws.Cells.Item(20, 5).Value = '  -0.25%  '
$ws.Cells.Item(21, 4).Value = '''0.987'
$ws.Cells.Item(21, 5).Value = '  -1.64%  '
$ws.Cells.Item(22, 4).Value = '''17.18'
$ws.Cells.Item(22, 5).Value = '  -2.58%  '
$ws.Cells.Item(23, 4).Value = '''106.78'
$ws.Cells.Item(23, 5).Value = '  +13.04%  '
$ws.Cells.Item(24, 4).Value = '''4.71'
$ws.Cells.Item(24, 5).Value = '  +3.00%  '
$ws.Cells.Item(25, 4).Value = '''5.07'
$ws.Cells.Item(25, 5).Value = '  +3.45%  '
$ws.Cells.Item(26, 4).Value = '''3.07'
$ws.Cells.Item(26, 5).Value = '  +4.60%  '
$ws.Cells.Item(27, 4).Value = '''11.00'
$ws.Cells.Item(27, 5).Value = '  -1.76%  '
$ws.Cells.Item(28, 4).Value = '''9.69'
$ws.Cells.Item(28, 5).Value = '  +4.32%  '
$ws.Cells.Item(29, 4).Value = '''33.87'
$ws.Cells.Item(29, 5).Value = '  +3.60%  '
$ws.Cells.Item(30, 5).Value = '  -3.06%  '
$ws.Cells.Item(31, 4).Value = '''12.50'
$ws.Cells.Item(31, 5).Value = '  +1.58%  '
$ws.Cells.Item(32, 5).Value = '  +2.87%  '
$ws.Cells.Item(33, 5).Value = '  +0.04%  '
$ws.Cells.Item(34, 4).Value = '''63.44'
$ws.Cells.Item(35, 4).Value = '''3.12'
$ws.Cells.Item(35, 5).Value = '  -5.77%  '
$ws.Cells.Item(37, 4).Value = '3.652.65'
$ws.Cells.Item(37, 5).Value = '  +0.77%  '
$ws.Cells.Item(38, 2).Value = 'Bittensor'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(38, 4).Value = '''513.82'
$ws.Cells.Item(38, 5).Value = '  -3.22%  '
$ws.Cells.Item(39, 4).Value = '''0.394'
$ws.Cells.Item(39, 5).Value = '  -4.09%  '
$ws.Cells.Item(40, 2).Value = 'Stacks'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(40, 4).Value = '''3.61'
$ws.Cells.Item(40, 5).Value = '  +5.10%  '
$ws.Cells.Item(41, 2).Value = 'PEPE'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(41, 4).Value = '0.0₃0786'
$ws.Cells.Item(41, 5).Value = '  -2.12%  '
$ws.Cells.Item(42, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(42, 4).Value = '''36.72'
$ws.Cells.Item(42, 5).Value = '  -4.36%  '
$ws.Cells.Item(43, 5).Value = '  +0.08%  '
$ws.Cells.Item(44, 5).Value = '  -1.06%  '
$ws.Cells.Item(45, 4).Value = '''2.92'
$ws.Cells.Item(45, 5).Value = '  -0.45%  '
$ws.Cells.Item(46, 4).Value = '''0.143'
$ws.Cells.Item(46, 5).Value = '  +3.34%  '
$ws.Cells.Item(47, 5).Value = '  -4.06%  '
$ws.Cells.Item(48, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(48, 4).Value = '''1.00'
$ws.Cells.Item(48, 5).Value = '  +0.40%  '
$ws.Cells.Item(49, 2).Value = 'THORChain'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(49, 4).Value = '''8.74'
$ws.Cells.Item(49, 5).Value = '  -5.55%  '
$ws.Cells.Item(50, 4).Value = '''131.92'
$ws.Cells.Item(50, 5).Value = '  -2.89%  '
$ws.Cells.Item(51, 4).Value = '''1.35'
$ws.Cells.Item(51, 5).Value = '  -6.63%  '
